$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = 7295
}

for ($r = 11; $r -le 252; $r++) {
    $ws.Cells.Item($r, 3).Value = 7293
}
